# Seminar 4 (scope & memory) update:
# On slide 3 ("Ways to define names"), the table "Table 3" has a row for
# "Function argument" whose example cell reads "Def func(value_1, value2):".
# Fix the capitalisation typo so it matches the Python keyword "def".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item("Table 3")  # "Table 3" graphic frame
$tbl = $shape.Table

$cell = $tbl.Cell(5, 2)             # row "Function argument", column "Example"
$tr = $cell.Shape.TextFrame.TextRange

# The cell's text is split across three runs: "Def ", "func", "(value_1, value2):".
# Re-assign the leading run's text to correct "Def " -> "def " while leaving the
# other runs (and their formatting) untouched.
$tr.Text = "def "
